# Normalize the "Recorded By" (column G) entries on the session analysis
# sheet: whenever "System" (or "admin@admin.com") was listed before the
# human editor, swap the order so the human editor's email comes first,
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "

        if ($parts.Length -eq 2) {
            $shouldSwap = $false

            if ($parts[0] -eq "System") {
                $shouldSwap = $true
            }
            if ($val -eq "admin@admin.com, dnasr281@gmail.com") {
                $shouldSwap = $true
            }

            if ($shouldSwap) {
                $cell.Value = $parts[1] + ", " + $parts[0]
            }
        }
    }
}
